$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 3.491328954696655
$ws.Range("B1").Value = 2.822103261947632
$ws.Range("C1").Value = 2.481647729873657
$ws.Range("D1").Value = 2.699149131774902
$ws.Range("E1").Value = 2.985780715942383
